$wb = $excel.ActiveWorkbook

# --- Append ", @deprecated" to a handful of scenario-outline Tags cells ---
$ws = $wb.Worksheets.Item("CUSTOMERWORKFLOWEU")
$ws.Range("C2").Value = "@FR, @FR031, @UR029, @Admin, @DesktopOnly, @deprecated"

$ws = $wb.Worksheets.Item("CUSTOMERWORKFLOWUK")
$ws.Range("C2").Value = "@FR, @FR030, @UR028, @Admin, @DesktopOnly, @deprecated"

$ws = $wb.Worksheets.Item("DESKTOPLICENSING")
$ws.Range("C2").Value = "@FR, @FR027, @UR025, @Admin, @DesktopOnly, @deprecated"

$ws = $wb.Worksheets.Item("SHELLINTEGRATION")
$ws.Range("C20").Value = "@FR035-03, @UR033-01, @deprecated"

$ws = $wb.Worksheets.Item("USERACCOUNTLICENSING")
$ws.Range("C14").Value = "@FR043-02, @UR041-01, @deprecated"

# --- Fix typos on WEBSITEADMINUSER ---
$ws = $wb.Worksheets.Item("WEBSITEADMINUSER")
$ws.Range("D9").Value = "I view the user roles"
$ws.Range("D10").Value = "there is a subscriber role for users that have registered but not setup an account"
$ws.Range("D39").Value = "I have created a CD Admin or CD User"

# --- Remove the stray duplicate "LoginEnabled" sub-header (+ its blank separator row)
#     that was left behind inside the NoLicenceForInvalidUser scenario on
#     USERACCOUNTLICENSING, shifting everything below it up by three rows. ---
$ws = $wb.Worksheets.Item("USERACCOUNTLICENSING")
$ws.Rows("24:26").Delete()
